$wb = $excel.ActiveWorkbook

$wsAbout   = $wb.Worksheets.Item("About")
$wsBts     = $wb.Worksheets.Item("BTS NTS Modal Profile Data")
$wsPass    = $wb.Worksheets.Item("AVLo-passengers")
$wsFreight = $wb.Worksheets.Item("AVLo-freight")

# The "AVLo-passengers" row that pulled from the BTS sheet's weighted
# (divide-by-10) helper row is being restored to point at the original
# un-normalized value row, before that helper row is removed below.
$wsPass.Cells.Item(5, 2).Formula = "='BTS NTS Modal Profile Data'!B36"

# Remove the extra "weighted value, adjusted for number of train cars per
# locomotive" row that was inserted on the BTS sheet. Deleting it shifts
# every row below up by one and Excel automatically repoints the other
# formulas that reference rows below it (e.g. AVLo-passengers!B7 and
# AVLo-freight!B6), and drops the now-unused shared string.
$wsBts.Rows("37").Delete()

# Reset the view state: clear the lingering scroll position/selection on
# the BTS sheet left over from editing it ...
$wsBts.Activate()
$wsBts.Range("A1").Select()

# ... and on the AVLo-passengers sheet too, which also loses the "active
# tab" flag.
$wsPass.Activate()
$wsPass.Range("A1").Select()

# "About" becomes the active sheet/tab again.
$wsAbout.Activate()
